$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.757.90'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.913.80'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '355.53'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.05'
$ws.Range('E6').Value = '  -3.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.561'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.97'
$ws.Range('E10').Value = '  -3.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0869'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.48'
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.74'
$ws.Range('E14').Value = '  -1.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.371.03'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.913.23'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.978'
$ws.Range('E17').Value = '  -2.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.687.41'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.51'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0977'
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.58'
$ws.Range('E23').Value = '  -0.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.27'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.79'
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.184'
$ws.Range('E26').Value = '  +11.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.83'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.55'
$ws.Range('E28').Value = '  +15.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +9.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.48'
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('E32').Value = '  -1.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.05'
$ws.Range('E33').Value = '  -1.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.15'
$ws.Range('E34').Value = '  -4.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '52.23'
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0439'
$ws.Range('E36').Value = '  -3.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.20'
$ws.Range('E38').Value = '  -3.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.21'
$ws.Range('E39').Value = '  -3.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.00'
$ws.Range('E40').Value = '  -4.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.72'
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.67'
$ws.Range('E43').Value = '  -5.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.11'
$ws.Range('E44').Value = '  -2.24%  '
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.48'
$ws.Range('E46').Value = '  -5.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.45'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.120.13'
$ws.Range('E48').Value = '  -3.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.249'
$ws.Range('E49').Value = '  -4.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0333'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('E51').Value = '  -0.45%  '
